$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - reorder column labels
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Update data rows 2-6 with the new one-hot values (row 7 unchanged)
$data = @{
    2 = @(0, 0, 0, 0, 1, 0);
    3 = @(0, 0, 0, 0, 0, 1);
    4 = @(1, 0, 0, 0, 0, 0);
    5 = @(0, 1, 0, 0, 0, 0);
    6 = @(0, 0, 1, 0, 0, 0);
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
